$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 3687.4736
$ws.Range("I62").Value = 3132.5
$ws.Range("J62").Value = 4638.857
$ws.Range("K62").Value = 3132.5
$ws.Range("L62").Value = 4638.857
$ws.Range("M62").Value = -2508.5
$ws.Range("N62").Value = -5886.857
$ws.Range("H65").Value = 3687.4736
$ws.Range("I65").Value = 3132.5
$ws.Range("J65").Value = 4638.857
$ws.Range("K65").Value = 15662.5
$ws.Range("L65").Value = 23194.285
$ws.Range("M65").Value = -12542.5
$ws.Range("N65").Value = -29434.285
$ws.Range("H129").Value = 556951.25
$ws.Range("J129").Value = 715995.1
$ws.Range("L129").Value = 2147985.3
$ws.Range("N129").Value = -2157985.3
$ws.Range("H132").Value = 2373.5278
$ws.Range("I132").Value = 2501.6765
$ws.Range("J132").Value = 195
$ws.Range("K132").Value = 7505.029500000001
$ws.Range("L132").Value = 585
$ws.Range("M132").Value = -4975.029500000001
$ws.Range("N132").Value = -5645
$ws.Range("H137").Value = 1943.6666
$ws.Range("I137").Value = 1679.5834
$ws.Range("J137").Value = 2207.75
$ws.Range("K137").Value = 5038.7502
$ws.Range("L137").Value = 6623.25
$ws.Range("M137").Value = -2488.7502
$ws.Range("N137").Value = -11723.25
$ws.Range("H138").Value = 2047.875
$ws.Range("I138").Value = 1515.3636
$ws.Range("J138").Value = 2249.862
$ws.Range("K138").Value = 4546.0908
$ws.Range("L138").Value = 6749.586
$ws.Range("M138").Value = 593.9092000000001
$ws.Range("N138").Value = -17029.586

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 100
$ws.Range("I4").Value = 100
$ws.Range("K4").Value = 100
$ws.Range("M4").Value = 16
$ws.Range("H5").Value = 290
$ws.Range("I5").Value = 290
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 290
$ws.Range("L5").Value = 0
$ws.Range("M5").ClearContents()
$ws.Range("N5").Value = -178
$ws.Range("H32").Value = 4796.6826
$ws.Range("I32").Value = 4671.035
$ws.Range("K32").Value = 4671.035
$ws.Range("M32").Value = -4384.035
$ws.Range("H45").Value = 3040.875
$ws.Range("I45").Value = 2792.611
$ws.Range("K45").Value = 2792.611
$ws.Range("M45").Value = -2415.611
$ws.Range("H96").Value = 0
$ws.Range("J96").Value = 0
$ws.Range("L96").ClearContents()
$ws.Range("N96").Value = 0
$ws.Range("H122").Value = 1428.826
$ws.Range("I122").Value = 1153.1177
$ws.Range("J122").Value = 2210
$ws.Range("K122").Value = 3459.3531
$ws.Range("L122").Value = 6630
$ws.Range("M122").Value = -1009.3531
$ws.Range("N122").Value = -11530
$ws.Range("H132").Value = 26755.184
$ws.Range("I132").Value = 1526.9788
$ws.Range("J132").Value = 117964.84
$ws.Range("K132").Value = 4580.936400000001
$ws.Range("L132").Value = 353894.52
$ws.Range("M132").Value = -2050.936400000001
$ws.Range("N132").Value = -358954.52

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 290
$ws.Range("I4").Value = 290
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 290
$ws.Range("L4").Value = 0
$ws.Range("M4").ClearContents()
$ws.Range("N4").Value = -175
$ws.Range("H80").Value = 515.5357
$ws.Range("J80").Value = 588.17645
$ws.Range("L80").Value = 588.17645
$ws.Range("N80").Value = -2584.17645
$ws.Range("H83").Value = 515.5357
$ws.Range("J83").Value = 588.17645
$ws.Range("L83").Value = 2940.88225
$ws.Range("N83").Value = -12924.88225
$ws.Range("H105").Value = 4411
$ws.Range("I105").Value = 5135
$ws.Range("K105").Value = 5135
$ws.Range("M105").Value = -3388
$ws.Range("H134").Value = 5654.7144
$ws.Range("I134").Value = 6325.5884
$ws.Range("J134").Value = 2803.5
$ws.Range("K134").Value = 18976.7652
$ws.Range("L134").Value = 8410.5
$ws.Range("M134").Value = -16441.7652
$ws.Range("N134").Value = -13480.5

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 56.666668
$ws.Range("I7").Value = 35
$ws.Range("K7").Value = 35
$ws.Range("M7").Value = 78
$ws.Range("H22").Value = 391
$ws.Range("I22").Value = 323.33334
$ws.Range("K22").Value = 323.33334
$ws.Range("M22").Value = 26.66665999999998
$ws.Range("H31").Value = 11996.925
$ws.Range("I31").Value = 14800.414
$ws.Range("J31").Value = 4605.909
$ws.Range("K31").Value = 14800.414
$ws.Range("L31").Value = 4605.909
$ws.Range("M31").Value = -14505.414
$ws.Range("N31").Value = -5195.909
$ws.Range("H34").Value = 11996.925
$ws.Range("I34").Value = 14800.414
$ws.Range("J34").Value = 4605.909
$ws.Range("K34").Value = 14800.414
$ws.Range("L34").Value = 4605.909
$ws.Range("M34").Value = -14598.414
$ws.Range("N34").Value = -5009.909
$ws.Range("H58").Value = 17090.42
$ws.Range("I58").Value = 1098.1305
$ws.Range("K58").Value = 1098.1305
$ws.Range("M58").Value = -895.1305
$ws.Range("H132").Value = 14987.282
$ws.Range("I132").Value = 16102.294
$ws.Range("K132").Value = 48306.882
$ws.Range("M132").Value = -45776.882
$ws.Range("H136").Value = 17090.42
$ws.Range("I136").Value = 1098.1305
$ws.Range("K136").Value = 3294.3915
$ws.Range("M136").Value = -744.3914999999997
$ws.Range("H141").Value = 0
$ws.Range("J141").Value = 0
$ws.Range("L141").ClearContents()
$ws.Range("N141").Value = 0

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H122").Value = 753
$ws.Range("I122").Value = 356.42856
$ws.Range("J122").Value = 916.2941
$ws.Range("K122").Value = 3207.85704
$ws.Range("L122").Value = 8246.6469
$ws.Range("M122").Value = -757.8570399999999
$ws.Range("N122").Value = -13146.6469
$ws.Range("H131").Value = 121297.484
$ws.Range("J131").Value = 127414.695
$ws.Range("L131").Value = 382244.085
$ws.Range("N131").Value = -392324.085

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 16197
$ws.Range("I70").Value = 34896
$ws.Range("J70").Value = 4977.6
$ws.Range("K70").Value = 34896
$ws.Range("L70").Value = 4977.6
$ws.Range("M70").Value = -34626
$ws.Range("N70").Value = -5517.6
$ws.Range("H73").Value = 16197
$ws.Range("I73").Value = 34896
$ws.Range("J73").Value = 4977.6
$ws.Range("K73").Value = 34896
$ws.Range("L73").Value = 4977.6
$ws.Range("M73").Value = -33960
$ws.Range("N73").Value = -6849.6
$ws.Range("H113").Value = 2172.6667
$ws.Range("I113").Value = 1539.2858
$ws.Range("K113").Value = 1539.2858
$ws.Range("M113").Value = 630.7141999999999
$ws.Range("H126").Value = 5101.1714
$ws.Range("I126").Value = 4897.696
$ws.Range("J126").Value = 5491.1665
$ws.Range("K126").Value = 14693.088
$ws.Range("L126").Value = 16473.4995
$ws.Range("M126").Value = -12223.088
$ws.Range("N126").Value = -21413.4995

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 1723.1111
$ws.Range("I93").Value = 1418
$ws.Range("J93").Value = 2333.3333
$ws.Range("K93").Value = 1418
$ws.Range("L93").Value = 2333.3333
$ws.Range("M93").Value = -170
$ws.Range("N93").Value = -4829.3333
$ws.Range("H136").Value = 14342.815
$ws.Range("I136").Value = 20960.68
$ws.Range("K136").Value = 62882.04
$ws.Range("M136").Value = -60332.04

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 1296.9
$ws.Range("I122").Value = 1093.9
$ws.Range("K122").Value = 3281.7
$ws.Range("M122").Value = -831.7000000000003
$ws.Range("H136").Value = 20834446
$ws.Range("I136").Value = 24391256
$ws.Range("J136").Value = 1693.5714
$ws.Range("K136").Value = 73173768
$ws.Range("L136").Value = 5080.7142
$ws.Range("M136").Value = -73171218
$ws.Range("N136").Value = -10180.7142
